$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "2025-11-03 Monday"; New = "2025-11-04 Tuesday" },
    @{ Old = "740×4=2960"; New = "714×3=2142" },
    @{ Old = "733×4=2932"; New = "575×3=1725" },
    @{ Old = "853×7=5971"; New = "137×2=274" },
    @{ Old = "871×9=7839"; New = "375×4=1500" },
    @{ Old = "321×4=1284"; New = "456×5=2280" },
    @{ Old = "553×3=1659"; New = "700×9=6300" },
    @{ Old = "803×5=4015"; New = "180×7=1260" },
    @{ Old = "965×6=5790"; New = "155×6=930" },
    @{ Old = "361×3=1083"; New = "868×7=6076" },
    @{ Old = "126×8=1008"; New = "854×6=5124" },
    @{ Old = "996×4=3984"; New = "110×8=880" },
    @{ Old = "418×4=1672"; New = "649×8=5192" },
    @{ Old = "710×6=4260"; New = "481×3=1443" },
    @{ Old = "166×9=1494"; New = "447×6=2682" },
    @{ Old = "486×4=1944"; New = "525×4=2100" },
    @{ Old = "654×5=3270"; New = "274×5=1370" },
    @{ Old = "538×6=3228"; New = "305×7=2135" },
    @{ Old = "486×9=4374"; New = "187×3=561" },
    @{ Old = "695×5=3475"; New = "967×7=6769" },
    @{ Old = "147×6=882"; New = "778×4=3112" },
    @{ Old = "611×3=1833"; New = "774×7=5418" },
    @{ Old = "678×2=1356"; New = "422×8=3376" },
    @{ Old = "435×6=2610"; New = "926×6=5556" },
    @{ Old = "483×7=3381"; New = "427×8=3416" },
    @{ Old = "186×3=558"; New = "910×4=3640" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
